$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record (dated 2021-09-09 / serial 44448) needs to be inserted
# right before the existing row 152 (dated 2020-12-02 / serial 44167), so
# insert a fresh row at position 152; this shifts the old rows 152-154 down
# to 153-155, which matches the target diff exactly (dimension grows from
# A1:R154 to A1:R155).
$ws.Rows(152).Insert()

# Populate the newly inserted row 152 with the new weekly data point.
$ws.Range("A152").Value = 7
$ws.Range("B152").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C152").Value = 'Ñuble'
$ws.Range("D152").Value = 44448
$ws.Range("E152").Value = 16
$ws.Range("F152").Value = 100112009
$ws.Range("G152").Value = 'Acelga'
$ws.Range("H152").Value = 'Sin especificar'
$ws.Range("I152").Value = 'Primera'
$ws.Range("J152").Value = 160
$ws.Range("K152").Value = 400
$ws.Range("L152").Value = 450
$ws.Range("M152").Value = 425
$ws.Range("N152").Value = '$/atado 1 a 1,5 kilos'
$ws.Range("O152").Value = 'Provincia de Diguillín'
$ws.Range("P152").Value = 283
$ws.Range("Q152").Value = 1.5
$ws.Range("R152").Value = 'Hortaliza'
